# Update "想去人数" (want-to-go count) figures across sheets to match the
# latest scrape output (gh-pages regeneration at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 26136
$ws1.Range("F7").Value = 585
$ws1.Range("F9").Value = 422
$ws1.Range("F11").Value = 347
$ws1.Range("F12").Value = 213
$ws1.Range("F16").Value = 38
$ws1.Range("F17").Value = 356
$ws1.Range("F18").Value = 50
$ws1.Range("F22").Value = 420

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 130
$ws2.Range("F8").Value = 106
$ws2.Range("F9").Value = 106
$ws2.Range("F10").Value = 428
$ws2.Range("F15").Value = 48

# 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 4947
$ws3.Range("F3").Value = 195

# 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 4947
$ws4.Range("F5").Value = 195
$ws4.Range("F6").Value = 26137
$ws4.Range("F11").Value = 585
$ws4.Range("F15").Value = 130
$ws4.Range("F16").Value = 130
$ws4.Range("F18").Value = 106
$ws4.Range("F19").Value = 106
$ws4.Range("F20").Value = 428
$ws4.Range("F21").Value = 422
$ws4.Range("F24").Value = 347
$ws4.Range("F25").Value = 213
$ws4.Range("F30").Value = 38
$ws4.Range("F33").Value = 356
$ws4.Range("F34").Value = 50
$ws4.Range("F35").Value = 48
$ws4.Range("F40").Value = 420
